$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09476766666666665
$ws.Range("I2").Value = 0.02336090049363864
$ws.Range("J2").Value = 0.02336090049363864
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3659943333333333
$ws.Range("N2").Value = 1.097983
$ws.Range("O2").Value = 0.006726051721149161
$ws.Range("P2").Value = 0.006726051721149162
$ws.Range("Q2").Value = 0.03468442898322222
$ws.Range("R2").Value = 0.3121598608489999
$ws.Range("S2").Value = 0.0001571266249728325
$ws.Range("T2").Value = 0.0001571266249728325
$ws.Range("G3").Value = 0.09476766666666665
$ws.Range("I3").Value = 0.02336090049363864
$ws.Range("J3").Value = 0.02336090049363864
$ws.Range("O3").Value = 0.001678071748088335
$ws.Range("P3").Value = 0.001678071748088335
$ws.Range("Q3").Value = 0.008653362000222221
$ws.Range("R3").Value = 0.07788025800199999
$ws.Range("S3").Value = 0.00003920126712827784
$ws.Range("T3").Value = 0.00003920126712827784
$ws.Range("G4").Value = 0.09476766666666665
$ws.Range("I4").Value = 0.02336090049363864
$ws.Range("J4").Value = 0.02336090049363864
$ws.Range("M4").Value = 53.897087
$ws.Range("N4").Value = 161.691261
$ws.Range("O4").Value = 0.9904923704135933
$ws.Range("P4").Value = 0.9904923704135934
$ws.Range("Q4").Value = 5.107701175120332
$ws.Range("R4").Value = 45.969310576083
$ws.Range("S4").Value = 0.02313879370494022
$ws.Range("T4").Value = 0.02313879370494022
$ws.Range("G5").Value = 0.09476766666666665
$ws.Range("I5").Value = 0.02336090049363864
$ws.Range("J5").Value = 0.02336090049363864
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06004666666666667
$ws.Range("N5").Value = 0.18014
$ws.Range("O5").Value = 0.001103506117169219
$ws.Range("P5").Value = 0.001103506117169219
$ws.Range("Q5").Value = 0.005690482491111111
$ws.Range("R5").Value = 0.05121434242
$ws.Range("S5").Value = 0.00002577889659731165
$ws.Range("T5").Value = 0.00002577889659731166
$ws.Range("H6").Value = 6.562189
$ws.Range("I6").Value = 0.5392086761288135
$ws.Range("J6").Value = 0.5392086761288134
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3659943333333333
$ws.Range("N6").Value = 1.097983
$ws.Range("O6").Value = 0.006726051721149161
$ws.Range("P6").Value = 0.006726051721149162
$ws.Range("Q6").Value = 0.8005746627541112
$ws.Range("R6").Value = 7.205171964787
$ws.Range("S6").Value = 0.003626745444134767
$ws.Range("T6").Value = 0.003626745444134767
$ws.Range("H7").Value = 6.562189
$ws.Range("I7").Value = 0.5392086761288135
$ws.Range("J7").Value = 0.5392086761288134
$ws.Range("O7").Value = 0.001678071748088335
$ws.Range("P7").Value = 0.001678071748088335
$ws.Range("S7").Value = 0.0009048308457358751
$ws.Range("T7").Value = 0.000904830845735875
$ws.Range("H8").Value = 6.562189
$ws.Range("I8").Value = 0.5392086761288135
$ws.Range("J8").Value = 0.5392086761288134
$ws.Range("M8").Value = 53.897087
$ws.Range("N8").Value = 161.691261
$ws.Range("O8").Value = 0.9904923704135933
$ws.Range("P8").Value = 0.9904923704135934
$ws.Range("Q8").Value = 117.8942904811477
$ws.Range("R8").Value = 1061.048614330329
$ws.Range("S8").Value = 0.534082079766404
$ws.Range("T8").Value = 0.534082079766404
$ws.Range("H9").Value = 6.562189
$ws.Range("I9").Value = 0.5392086761288135
$ws.Range("J9").Value = 0.5392086761288134
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.06004666666666667
$ws.Range("N9").Value = 0.18014
$ws.Range("O9").Value = 0.001103506117169219
$ws.Range("P9").Value = 0.001103506117169219
$ws.Range("Q9").Value = 0.1313458584955556
$ws.Range("R9").Value = 1.18211272646
$ws.Range("S9").Value = 0.0005950200725388617
$ws.Range("T9").Value = 0.0005950200725388617
$ws.Range("G10").Value = 1.774514666666667
$ws.Range("H10").Value = 5.323544
$ws.Range("I10").Value = 0.437430423377548
$ws.Range("J10").Value = 0.437430423377548
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3659943333333333
$ws.Range("N10").Value = 1.097983
$ws.Range("O10").Value = 0.006726051721149161
$ws.Range("P10").Value = 0.006726051721149162
$ws.Range("Q10").Value = 0.6494623124168889
$ws.Range("R10").Value = 5.845160811752
$ws.Range("S10").Value = 0.002942179652041563
$ws.Range("T10").Value = 0.002942179652041563
$ws.Range("G11").Value = 1.774514666666667
$ws.Range("H11").Value = 5.323544
$ws.Range("I11").Value = 0.437430423377548
$ws.Range("J11").Value = 0.437430423377548
$ws.Range("O11").Value = 0.001678071748088335
$ws.Range("P11").Value = 0.001678071748088335
$ws.Range("Q11").Value = 0.1620333002328889
$ws.Range("R11").Value = 1.458299702096
$ws.Range("S11").Value = 0.0007340396352241825
$ws.Range("T11").Value = 0.0007340396352241825
$ws.Range("G12").Value = 1.774514666666667
$ws.Range("H12").Value = 5.323544
$ws.Range("I12").Value = 0.437430423377548
$ws.Range("J12").Value = 0.437430423377548
$ws.Range("M12").Value = 53.897087
$ws.Range("N12").Value = 161.691261
$ws.Range("O12").Value = 0.9904923704135933
$ws.Range("P12").Value = 0.9904923704135934
$ws.Range("Q12").Value = 95.64117137210933
$ws.Range("R12").Value = 860.770542348984
$ws.Range("S12").Value = 0.4332714969422491
$ws.Range("T12").Value = 0.4332714969422492
$ws.Range("G13").Value = 1.774514666666667
$ws.Range("H13").Value = 5.323544
$ws.Range("I13").Value = 0.437430423377548
$ws.Range("J13").Value = 0.437430423377548
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.06004666666666667
$ws.Range("N13").Value = 0.18014
$ws.Range("O13").Value = 0.001103506117169219
$ws.Range("P13").Value = 0.001103506117169219
$ws.Range("Q13").Value = 0.1065536906844445
$ws.Range("R13").Value = 0.9589832161600002
$ws.Range("S13").Value = 0.0004827071480330453
$ws.Range("T13").Value = 0.0004827071480330454
